# Rotate the data rows (2-5) down by one, with the last row (5)
# wrapping around to become the new row 2:
#   new row2 = old row5
#   new row3 = old row2
#   new row4 = old row3
#   new row5 = old row4
#
# Only the columns whose values actually differ between the four rows
# need to be touched; columns that already hold the same value in every
# row (location, date, observer, etc.) are left completely alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 5

# Columns that are stored as plain numbers.
$numericCols = @("A", "B", "E", "Q", "R")

# Columns that are stored as text (some look numeric, e.g. "12", so we
# force a text format before writing to stop Excel from re-interpreting
# them as numbers/dates).
$textCols = @("D", "F", "G", "H", "I", "J", "K", "M", "N", "AF")

$allCols = $numericCols + $textCols

# 1. Snapshot the current values of the cells that are about to move.
$snapshot = @{}
foreach ($col in $allCols) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# 2. Work out, for each destination row, which row currently holds the
#    data that needs to end up there.
$sourceRowFor = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($r -eq $firstRow) {
        $sourceRowFor[$r] = $lastRow
    } else {
        $sourceRowFor[$r] = $r - 1
    }
}

# 3. Write the numeric columns back using their plain values.
foreach ($col in $numericCols) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $src = $sourceRowFor[$r]
        $ws.Range("$col$r").Value2 = $snapshot["$col$src"]
    }
}

# 4. Write the text columns back, forcing a text number format so values
#    such as "12" or empty strings stay text instead of becoming numbers.
foreach ($col in $textCols) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $src = $sourceRowFor[$r]
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value2 = $snapshot["$col$src"]
        $cell.NumberFormat = "General"
    }
}
